$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 500
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 471.3611
$ws.Range("I33").Value = 245.04167
$ws.Range("K33").Value = 245.04167
$ws.Range("M33").Value = -16.04167000000001
$ws.Range("H106").Value = 2782301.2
$ws.Range("I106").Value = 3513170
$ws.Range("K106").Value = 3513170
$ws.Range("M106").Value = -3512539
$ws.Range("H137").Value = 26801.727
$ws.Range("I137").Value = 45188.316
$ws.Range("J137").Value = 1848.5
$ws.Range("K137").Value = 135564.948
$ws.Range("L137").Value = 5545.5
$ws.Range("M137").Value = -133014.948
$ws.Range("N137").Value = -10645.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9340517
$ws.Range("I32").Value = 2904837.8
$ws.Range("K32").Value = 2904837.8
$ws.Range("M32").Value = -2904550.8
$ws.Range("H45").Value = 3078.4666
$ws.Range("I45").Value = 3053.389
$ws.Range("K45").Value = 3053.389
$ws.Range("M45").Value = -2676.389
$ws.Range("H61").Value = 3005.111
$ws.Range("I61").Value = 2818.25
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 2818.25
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -2606.25
$ws.Range("N61").Value = -4924
$ws.Range("H74").Value = 1997.25
$ws.Range("I74").Value = 1709.4375
$ws.Range("J74").Value = 4299.75
$ws.Range("K74").Value = 1709.4375
$ws.Range("L74").Value = 4299.75
$ws.Range("M74").Value = -835.4375
$ws.Range("N74").Value = -6047.75
$ws.Range("H77").Value = 1997.25
$ws.Range("I77").Value = 1709.4375
$ws.Range("J77").Value = 4299.75
$ws.Range("K77").Value = 8547.1875
$ws.Range("L77").Value = 21498.75
$ws.Range("M77").Value = -4179.1875
$ws.Range("N77").Value = -30234.75
$ws.Range("H97").Value = 864.1667
$ws.Range("I97").Value = 543.4
$ws.Range("J97").Value = 2468
$ws.Range("K97").Value = 543.4
$ws.Range("L97").Value = 2468
$ws.Range("M97").Value = -47.39999999999998
$ws.Range("N97").Value = -3460
$ws.Range("H136").Value = 3005.111
$ws.Range("I136").Value = 2818.25
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 8454.75
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -5904.75
$ws.Range("N136").Value = -18600
$ws.Range("H137").Value = 133555.22
$ws.Range("I137").Value = 75999.25
$ws.Range("K137").Value = 75999.25
$ws.Range("M137").Value = -70899.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 989.6383
$ws.Range("I94").Value = 919.08826
$ws.Range("K94").Value = 919.08826
$ws.Range("M94").Value = -468.08826
$ws.Range("H105").Value = 2749.0652
$ws.Range("I105").Value = 3240.56
$ws.Range("J105").Value = 2163.9524
$ws.Range("K105").Value = 3240.56
$ws.Range("L105").Value = 2163.9524
$ws.Range("M105").Value = -1493.56
$ws.Range("N105").Value = -5657.9524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 341.77777
$ws.Range("J22").Value = 221.5
$ws.Range("L22").Value = 221.5
$ws.Range("N22").Value = -921.5
$ws.Range("H31").Value = 3977.4878
$ws.Range("I31").Value = 2590.862
$ws.Range("K31").Value = 2590.862
$ws.Range("M31").Value = -2295.862
$ws.Range("H34").Value = 3977.4878
$ws.Range("I34").Value = 2590.862
$ws.Range("K34").Value = 2590.862
$ws.Range("M34").Value = -2388.862
$ws.Range("H58").Value = 2821.5476
$ws.Range("I58").Value = 2196.6206
$ws.Range("J58").Value = 4215.615
$ws.Range("K58").Value = 2196.6206
$ws.Range("L58").Value = 4215.615
$ws.Range("M58").Value = -1993.6206
$ws.Range("N58").Value = -4621.615
$ws.Range("H132").Value = 1830
$ws.Range("J132").Value = 1795
$ws.Range("L132").Value = 5385
$ws.Range("N132").Value = -10445
$ws.Range("H133").Value = 98900
$ws.Range("J133").Value = 98900
$ws.Range("L133").Value = 98900
$ws.Range("N133").Value = -103960
$ws.Range("H134").Value = 2034.0834
$ws.Range("J134").Value = 3705.7144
$ws.Range("L134").Value = 11117.1432
$ws.Range("N134").Value = -16187.1432
$ws.Range("H136").Value = 2821.5476
$ws.Range("I136").Value = 2196.6206
$ws.Range("J136").Value = 4215.615
$ws.Range("K136").Value = 6589.861800000001
$ws.Range("L136").Value = 12646.845
$ws.Range("M136").Value = -4039.861800000001
$ws.Range("N136").Value = -17746.845
$ws.Range("H137").Value = 85000.2
$ws.Range("J137").Value = 93750.25
$ws.Range("L137").Value = 93750.25
$ws.Range("N137").Value = -103950.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1875.75
$ws.Range("I68").Value = 1001.5
$ws.Range("J68").Value = 2750
$ws.Range("K68").Value = 3004.5
$ws.Range("L68").Value = 8250
$ws.Range("M68").Value = -2193.5
$ws.Range("N68").Value = -9872
$ws.Range("H71").Value = 1875.75
$ws.Range("I71").Value = 1001.5
$ws.Range("J71").Value = 2750
$ws.Range("K71").Value = 9013.5
$ws.Range("L71").Value = 24750
$ws.Range("M71").Value = -4957.5
$ws.Range("N71").Value = -32862

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 311.60715
$ws.Range("I2").Value = 181.22223
$ws.Range("J2").Value = 546.3
$ws.Range("K2").Value = 181.22223
$ws.Range("L2").Value = 546.3
$ws.Range("M2").Value = -68.22223
$ws.Range("N2").Value = -772.3
$ws.Range("H119").Value = 92250
$ws.Range("J119").Value = 69666.664
$ws.Range("L119").Value = 69666.664
$ws.Range("N119").Value = -79342.664
$ws.Range("H123").Value = 29258.666
$ws.Range("J123").Value = 29258.666
$ws.Range("L123").Value = 29258.666
$ws.Range("N123").Value = -34158.666
$ws.Range("H126").Value = 1732.1765
$ws.Range("I126").Value = 1667.2
$ws.Range("J126").Value = 1825
$ws.Range("K126").Value = 5001.6
$ws.Range("L126").Value = 5475
$ws.Range("M126").Value = -2531.6
$ws.Range("N126").Value = -10415

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 4953
$ws.Range("I11").Value = 4953
$ws.Range("K11").Value = 4953
$ws.Range("M11").Value = -4813
$ws.Range("H36").Value = 92703
$ws.Range("J36").Value = 92703
$ws.Range("L36").Value = 92703
$ws.Range("N36").Value = -93827
$ws.Range("H40").Value = 4300.533
$ws.Range("I40").Value = 4428.68
$ws.Range("J40").Value = 3659.8
$ws.Range("K40").Value = 4428.68
$ws.Range("L40").Value = 3659.8
$ws.Range("M40").Value = -4292.68
$ws.Range("N40").Value = -3931.8
$ws.Range("H99").Value = 70134
$ws.Range("I99").Value = 42499.75
$ws.Range("K99").Value = 42499.75
$ws.Range("M99").Value = -39504.75
$ws.Range("H111").Value = 85984.664
$ws.Range("J111").Value = 85984.664
$ws.Range("L111").Value = 85984.664
$ws.Range("N111").Value = -94164.664
$ws.Range("H132").Value = 3888.5593
$ws.Range("I132").Value = 2899.875
$ws.Range("K132").Value = 8699.625
$ws.Range("M132").Value = -6169.625
$ws.Range("H136").Value = 5450.12
$ws.Range("I136").Value = 4977
$ws.Range("K136").Value = 14931
$ws.Range("M136").Value = -12381
$ws.Range("H141").Value = 200000
$ws.Range("J141").Value = 200000
$ws.Range("L141").Value = 200000
$ws.Range("N141").Value = -210360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 15606.4
$ws.Range("I20").Value = 6003.3335
$ws.Range("K20").Value = 6003.3335
$ws.Range("M20").Value = -5763.3335
$ws.Range("H122").Value = 1899.871
$ws.Range("I122").Value = 1828.75
$ws.Range("K122").Value = 5486.25
$ws.Range("M122").Value = -3036.25
$ws.Range("H128").Value = 200000
$ws.Range("J128").Value = 200000
$ws.Range("L128").Value = 200000
$ws.Range("N128").Value = -209960
$ws.Range("H136").Value = 24571.49
$ws.Range("I136").Value = 1925.15
$ws.Range("K136").Value = 5775.450000000001
$ws.Range("M136").Value = -3225.450000000001
